$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increase barrel multiplier (J2) from 1.5 to 2.5; dependent formulas recalc automatically
$ws.Range("J2").Value = 2.5

# Update the selected/active cell to reflect where the editor left off
$ws.Range("H11").Select()
